# feat: make exports german for logbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header labels
$ws.Range("A2").Value = "Report Typ:"
$ws.Range("B2").Value = "Fahrtenbuch"
$ws.Range("A5").Value = "Gruppe"
$ws.Range("A6").Value = "Periode:"
$ws.Range("A7").Value = "Einträge:"

# Summary labels
$ws.Range("A8").Value = "Distanz total:"
$ws.Range("A9").Value = "Distanz privat"
$ws.Range("A10").Value = "Distanz geschäftlich"
$ws.Range("D8").Value = "Dauer total"
$ws.Range("D9").Value = "Dauer privat"
$ws.Range("D10").Value = "Dauer geschäftlich"

# Table column headers (row 12)
$ws.Range("B12").Value = "Start Adresse"
$ws.Range("C12").Value = "Kilometerstand Start"
$ws.Range("D12").Value = "Ende"
$ws.Range("E12").Value = "Ende Adresse"
$ws.Range("F12").Value = "Kilometerstand Ende"
$ws.Range("G12").Value = "Dauer"
$ws.Range("H12").Value = "Distanz"
$ws.Range("I12").Value = "Höchstgeschwindigkeit"
$ws.Range("J12").Value = "Durchschnittsgeschwindigkeit"
$ws.Range("K12").Value = "Kraftstoffverbrauch"
$ws.Range("L12").Value = "Fahrer"
$ws.Range("M12").Value = "Typ"
$ws.Range("N12").Value = "Notiz"

# Restore the cell selection seen in the saved workbook
$ws.Range("D9").Select()
